$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = -2378

$ws.Range("D6").Value = 110
$ws.Range("E6").Value = 105

$ws.Range("D7").Value = 110
$ws.Range("E7").Value = 113

$ws.Range("D8").Value = 110
$ws.Range("E8").Value = 108

$ws.Range("D9").Value = 110
$ws.Range("E9").Value = 126

$ws.Range("D12").Value = 144
$ws.Range("E12").Value = 129

$ws.Range("D13").Value = 110
$ws.Range("E13").Value = 120

$ws.Range("D14").Value = 110
$ws.Range("E14").Value = 106

$ws.Range("D15").Value = 110
$ws.Range("E15").Value = 107

$ws.Range("D16").Value = 110
$ws.Range("E16").Value = 106

$ws.Range("D19").Value = 144
$ws.Range("E19").Value = 124

$ws.Range("D20").Value = 110
$ws.Range("E20").Value = 143

$ws.Range("D21").Value = 110
$ws.Range("E21").Value = 115

$ws.Range("D22").Value = 110
$ws.Range("E22").Value = 134

$ws.Range("D23").Value = 110
$ws.Range("E23").Value = 87

$ws.Range("D26").Value = 110
$ws.Range("D27").Value = 110
$ws.Range("D28").Value = 110
$ws.Range("D29").Value = 110
$ws.Range("D30").Value = 110

$ws.Range("D33").Value = 110
$ws.Range("D34").Value = 41
